$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Sentiment column (K) from "Negative" to "Neutral" for the rows
# that represent the newly-excluded "International" coverage, per the
# "Added International Brand Exclusion" change.
$rowsToUpdate = @(2, 3, 5, 6, 8, 9, 10, 11, 12, 13, 14)
foreach ($row in $rowsToUpdate) {
    $ws.Range("K$row").Value = "Neutral"
}

# Move the active selection / view from AA8 (with R1 scrolled to the
# top-left) back to C9 with no special scroll position.
$ws.Range("C9").Select()
